$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D4").Value = "파이썬 딥러닝 텐서플로 책이 출간되었습니다!"
$ws.Range("E4").Value = "https://teddylee777.github.io/thoughts/tf-book"

$ws.Range("D28").Value = "Evolutionary strategy 1 - Simple Gaussian Evolution"
$ws.Range("E28").Value = "https://ropiens.tistory.com/137"

$ws.Range("D37").Value = "[paper review] Do Transformers Really Perform Bad for Graph Representation?"
$ws.Range("E37").Value = "http://dsba.korea.ac.kr/seminar/?uid=1785&mod=document&pageid=1"

$ws.Range("D50").Value = "사이언스/네이처, 같은날 단백질 접힘 프로토콜 각각 발표 RoseTTAFold / AlphaFold"
$ws.Range("E50").Value = "http://incredible.egloos.com/7521788"
